# Generate Report for Handoff
#
# Updates the localization-status workbook so that:
#  - the previously in-flight "174a2b1c...md" file is now reported as a
#    dependency-tracked "cce5b3b7...md" file (with two new ".png"
#    dependencies: "3dd3ed8b...png" and "460c6bab...png"),
#  - the Overview / zh-cn / de-de sheets each grow from one data row to
#    three data rows (one per tracked file),
#  - handoff datetimes move forward to the new report run.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "3dd3ed8b-b50e-4fbe-b2b5-e534b4117dd2.png"
$ov.Range("D2").Value = "2016-03-22 19:04:07"

$ov.Range("A3").Value = "460c6bab-9d4f-4c29-9985-77fe979108ef.png"
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-03-22 19:04:07"

$ov.Range("A4").Value = "cce5b3b7-6bcc-4759-ad91-d12e32281e35.md"
$ov.Range("B4").Value = "Ready for handoff"
$ov.Range("C4").Value = "Ready for handoff"
$ov.Range("D4").Value = "2016-03-22 19:04:07"

$ov.Hyperlinks.Item(1).Address = "https://github.com/OpenLocalizationTest/oltest/blob/b9493ce9af5be9c3813788c16640c13a952429a6/e2e/3dd3ed8b-b50e-4fbe-b2b5-e534b4117dd2.png"
$ov.Hyperlinks.Item(1).TextToDisplay = "3dd3ed8b-b50e-4fbe-b2b5-e534b4117dd2.png"

$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b9493ce9af5be9c3813788c16640c13a952429a6/e2e/460c6bab-9d4f-4c29-9985-77fe979108ef.png", "", "", "460c6bab-9d4f-4c29-9985-77fe979108ef.png")
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/b9493ce9af5be9c3813788c16640c13a952429a6/e2e/cce5b3b7-6bcc-4759-ad91-d12e32281e35.md", "", "", "cce5b3b7-6bcc-4759-ad91-d12e32281e35.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "3dd3ed8b-b50e-4fbe-b2b5-e534b4117dd2.png"
$zh.Range("B2").Value = ".png"
$zh.Range("D2").Value = "8821e1065099484c1a0ce2ee997aa204f21717be.png"
$zh.Range("E2").Value = "2016-03-22 19:04:01"
$zh.Range("J2").Value = "IsDependency"
$zh.Range("K2").Value = "e2e\cce5b3b7-6bcc-4759-ad91-d12e32281e35.md"

$zh.Range("A3").Value = "460c6bab-9d4f-4c29-9985-77fe979108ef.png"
$zh.Range("B3").Value = ".png"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("D3").Value = "d02dc8587e7c3bebdcd9e4f19d1e0aae36034504.png"
$zh.Range("E3").Value = "2016-03-22 19:04:01"
$zh.Range("H3").Value = "0001-01-01 00:00:00"
$zh.Range("J3").Value = "IsDependency"
$zh.Range("K3").Value = "e2e\cce5b3b7-6bcc-4759-ad91-d12e32281e35.md"

$zh.Range("A4").Value = "cce5b3b7-6bcc-4759-ad91-d12e32281e35.md"
$zh.Range("B4").Value = ".md"
$zh.Range("C4").Value = "Ready for handoff"
$zh.Range("D4").Value = "cce5b3b7-6bcc-4759-ad91-d12e32281e35.7161c2e5ff0d0712ed7800f840b7a07091cf4dfe.zh-cn.xlf"
$zh.Range("E4").Value = "2016-03-22 19:04:01"
$zh.Range("H4").Value = "0001-01-01 00:00:00"
$zh.Range("J4").Value = "Include"

$zh.Hyperlinks.Item(1).Address = "https://github.com/OpenLocalizationTest/oltest/blob/b9493ce9af5be9c3813788c16640c13a952429a6/e2e/3dd3ed8b-b50e-4fbe-b2b5-e534b4117dd2.png"
$zh.Hyperlinks.Item(1).TextToDisplay = "3dd3ed8b-b50e-4fbe-b2b5-e534b4117dd2.png"
$zh.Hyperlinks.Item(2).Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/35429f74f10731d62c0a9c9271fbb3ecde727ca5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/8821e1065099484c1a0ce2ee997aa204f21717be.png"
$zh.Hyperlinks.Item(2).TextToDisplay = "8821e1065099484c1a0ce2ee997aa204f21717be.png"

$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b9493ce9af5be9c3813788c16640c13a952429a6/e2e/460c6bab-9d4f-4c29-9985-77fe979108ef.png", "", "", "460c6bab-9d4f-4c29-9985-77fe979108ef.png")
$zh.Hyperlinks.Add($zh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/35429f74f10731d62c0a9c9271fbb3ecde727ca5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d02dc8587e7c3bebdcd9e4f19d1e0aae36034504.png", "", "", "d02dc8587e7c3bebdcd9e4f19d1e0aae36034504.png")
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/b9493ce9af5be9c3813788c16640c13a952429a6/e2e/cce5b3b7-6bcc-4759-ad91-d12e32281e35.md", "", "", "cce5b3b7-6bcc-4759-ad91-d12e32281e35.md")
$zh.Hyperlinks.Add($zh.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/35429f74f10731d62c0a9c9271fbb3ecde727ca5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/cce5b3b7-6bcc-4759-ad91-d12e32281e35.7161c2e5ff0d0712ed7800f840b7a07091cf4dfe.zh-cn.xlf", "", "", "cce5b3b7-6bcc-4759-ad91-d12e32281e35.7161c2e5ff0d0712ed7800f840b7a07091cf4dfe.zh-cn.xlf")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "3dd3ed8b-b50e-4fbe-b2b5-e534b4117dd2.png"
$de.Range("B2").Value = ".png"
$de.Range("D2").Value = "8821e1065099484c1a0ce2ee997aa204f21717be.png"
$de.Range("E2").Value = "2016-03-22 19:04:07"
$de.Range("J2").Value = "IsDependency"
$de.Range("K2").Value = "e2e\cce5b3b7-6bcc-4759-ad91-d12e32281e35.md"

$de.Range("A3").Value = "460c6bab-9d4f-4c29-9985-77fe979108ef.png"
$de.Range("B3").Value = ".png"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("D3").Value = "d02dc8587e7c3bebdcd9e4f19d1e0aae36034504.png"
$de.Range("E3").Value = "2016-03-22 19:04:07"
$de.Range("H3").Value = "0001-01-01 00:00:00"
$de.Range("J3").Value = "IsDependency"
$de.Range("K3").Value = "e2e\cce5b3b7-6bcc-4759-ad91-d12e32281e35.md"

$de.Range("A4").Value = "cce5b3b7-6bcc-4759-ad91-d12e32281e35.md"
$de.Range("B4").Value = ".md"
$de.Range("C4").Value = "Ready for handoff"
$de.Range("D4").Value = "cce5b3b7-6bcc-4759-ad91-d12e32281e35.7161c2e5ff0d0712ed7800f840b7a07091cf4dfe.de-de.xlf"
$de.Range("E4").Value = "2016-03-22 19:04:07"
$de.Range("H4").Value = "0001-01-01 00:00:00"
$de.Range("J4").Value = "Include"

$de.Hyperlinks.Item(1).Address = "https://github.com/OpenLocalizationTest/oltest/blob/b9493ce9af5be9c3813788c16640c13a952429a6/e2e/3dd3ed8b-b50e-4fbe-b2b5-e534b4117dd2.png"
$de.Hyperlinks.Item(1).TextToDisplay = "3dd3ed8b-b50e-4fbe-b2b5-e534b4117dd2.png"
$de.Hyperlinks.Item(2).Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/184e1fbfdc060fec9a92321cdd148f2c8613b987/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/8821e1065099484c1a0ce2ee997aa204f21717be.png"
$de.Hyperlinks.Item(2).TextToDisplay = "8821e1065099484c1a0ce2ee997aa204f21717be.png"

$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b9493ce9af5be9c3813788c16640c13a952429a6/e2e/460c6bab-9d4f-4c29-9985-77fe979108ef.png", "", "", "460c6bab-9d4f-4c29-9985-77fe979108ef.png")
$de.Hyperlinks.Add($de.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/184e1fbfdc060fec9a92321cdd148f2c8613b987/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d02dc8587e7c3bebdcd9e4f19d1e0aae36034504.png", "", "", "d02dc8587e7c3bebdcd9e4f19d1e0aae36034504.png")
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/b9493ce9af5be9c3813788c16640c13a952429a6/e2e/cce5b3b7-6bcc-4759-ad91-d12e32281e35.md", "", "", "cce5b3b7-6bcc-4759-ad91-d12e32281e35.md")
$de.Hyperlinks.Add($de.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/184e1fbfdc060fec9a92321cdd148f2c8613b987/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/cce5b3b7-6bcc-4759-ad91-d12e32281e35.7161c2e5ff0d0712ed7800f840b7a07091cf4dfe.de-de.xlf", "", "", "cce5b3b7-6bcc-4759-ad91-d12e32281e35.7161c2e5ff0d0712ed7800f840b7a07091cf4dfe.de-de.xlf")

Write-Host "Report regenerated for handoff"
